$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Find the last used row in column A (the "Beteckning" column) to know the data extent.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C ("Förändrad") holds a date value (Excel serial 45171 -> 2023-09-02) for every
# data row (rows 2..lastRow). The commit updates this "last changed" date to 45172
# (2023-09-03) for all rows.
$range = $ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item($lastRow, 3))
$range.Value = 45172
